$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 2949.5
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 2949.5
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 8848.5
$ws.Range("N70").Value = -9388.5
$ws.Range("M70").Value = $null
# Row 73
$ws.Range("H73").Value = 2949.5
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 2949.5
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 8848.5
$ws.Range("N73").Value = -10720.5
$ws.Range("M73").Value = $null
# Row 99
$ws.Range("H99").Value = 687.8333
$ws.Range("I99").Value = 650.4
$ws.Range("J99").Value = 875
$ws.Range("K99").Value = 1951.2
$ws.Range("L99").Value = 2625
$ws.Range("M99").Value = -453.1999999999998
$ws.Range("N99").Value = -5621
# Row 103
$ws.Range("H103").Value = 512.53845
$ws.Range("I103").Value = 263.6
$ws.Range("J103").Value = 668.125
$ws.Range("K103").Value = 790.8000000000001
$ws.Range("L103").Value = 2004.375
$ws.Range("M103").Value = -204.8000000000001
$ws.Range("N103").Value = -3176.375
# Row 132
$ws.Range("H132").Value = 5650.4
$ws.Range("I132").Value = 4214
$ws.Range("J132").Value = 9002
$ws.Range("K132").Value = 12642
$ws.Range("L132").Value = 27006
$ws.Range("M132").Value = -10112
$ws.Range("N132").Value = -32066

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 4227.9165
$ws.Range("I2").Value = 4260.75
$ws.Range("J2").Value = 4162.25
$ws.Range("K2").Value = 4260.75
$ws.Range("L2").Value = 4162.25
$ws.Range("M2").Value = -4147.75
$ws.Range("N2").Value = -4388.25
# Row 32
$ws.Range("H32").Value = 2957.7368
$ws.Range("I32").Value = 2217.4
$ws.Range("J32").Value = 11595
$ws.Range("K32").Value = 2217.4
$ws.Range("L32").Value = 11595
$ws.Range("M32").Value = -1930.4
$ws.Range("N32").Value = -12169
# Row 35
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = $null
# Row 61
$ws.Range("H61").Value = 2245.0476
$ws.Range("I61").Value = 1991.8948
$ws.Range("J61").Value = 4650
$ws.Range("K61").Value = 1991.8948
$ws.Range("L61").Value = 4650
$ws.Range("M61").Value = -1779.8948
$ws.Range("N61").Value = -5074
# Row 116
$ws.Range("H116").Value = 4227.9165
$ws.Range("I116").Value = 4260.75
$ws.Range("J116").Value = 4162.25
$ws.Range("K116").Value = 4260.75
$ws.Range("L116").Value = 4162.25
$ws.Range("M116").Value = -1966.75
$ws.Range("N116").Value = -8750.25
# Row 122
$ws.Range("H122").Value = 10615.5
$ws.Range("I122").Value = 8738.799999999999
$ws.Range("J122").Value = 19999
$ws.Range("K122").Value = 26216.4
$ws.Range("L122").Value = 59997
$ws.Range("M122").Value = -23766.4
$ws.Range("N122").Value = -64897
# Row 132
$ws.Range("H132").Value = 2043.2354
$ws.Range("I132").Value = 1860.375
$ws.Range("J132").Value = 4969
$ws.Range("K132").Value = 5581.125
$ws.Range("L132").Value = 14907
$ws.Range("M132").Value = -3051.125
$ws.Range("N132").Value = -19967
# Row 136
$ws.Range("H136").Value = 2245.0476
$ws.Range("I136").Value = 1991.8948
$ws.Range("J136").Value = 4650
$ws.Range("K136").Value = 5975.6844
$ws.Range("L136").Value = 13950
$ws.Range("M136").Value = -3425.6844
$ws.Range("N136").Value = -19050

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 4227.9165
$ws.Range("I3").Value = 4260.75
$ws.Range("J3").Value = 4162.25
$ws.Range("K3").Value = 4260.75
$ws.Range("L3").Value = 4162.25
$ws.Range("M3").Value = -4146.75
$ws.Range("N3").Value = -4390.25
# Row 20
$ws.Range("H20").Value = 4332.6665
$ws.Range("I20").Value = 4199.2
$ws.Range("J20").Value = 5000
$ws.Range("K20").Value = 4199.2
$ws.Range("L20").Value = 5000
$ws.Range("M20").Value = -3952.2
$ws.Range("N20").Value = -5494
# Row 80
$ws.Range("H80").Value = 499.06668
$ws.Range("I80").Value = 662.4
$ws.Range("J80").Value = 417.4
$ws.Range("K80").Value = 662.4
$ws.Range("L80").Value = 417.4
$ws.Range("M80").Value = 335.6
$ws.Range("N80").Value = -2413.4
# Row 82
$ws.Range("H82").Value = 30339
$ws.Range("I82").Value = 13473.75
$ws.Range("J82").Value = 97800
$ws.Range("K82").Value = 13473.75
$ws.Range("L82").Value = 97800
$ws.Range("M82").Value = -13090.75
$ws.Range("N82").Value = -98566
# Row 83
$ws.Range("H83").Value = 499.06668
$ws.Range("I83").Value = 662.4
$ws.Range("J83").Value = 417.4
$ws.Range("K83").Value = 3312
$ws.Range("L83").Value = 2087
$ws.Range("M83").Value = 1680
$ws.Range("N83").Value = -12071
# Row 85
$ws.Range("H85").Value = 30339
$ws.Range("I85").Value = 13473.75
$ws.Range("J85").Value = 97800
$ws.Range("K85").Value = 13473.75
$ws.Range("L85").Value = 97800
$ws.Range("M85").Value = -12147.75
$ws.Range("N85").Value = -100452
# Row 99
$ws.Range("H99").Value = 4772.067
$ws.Range("I99").Value = 4863.2856
$ws.Range("J99").Value = 3495
$ws.Range("K99").Value = 4863.2856
$ws.Range("L99").Value = 3495
$ws.Range("M99").Value = -3365.2856
$ws.Range("N99").Value = -6491

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 8356.429
$ws.Range("I62").Value = 9339
$ws.Range("J62").Value = 5900
$ws.Range("K62").Value = 9339
$ws.Range("L62").Value = 5900
$ws.Range("M62").Value = -8715
$ws.Range("N62").Value = -7148
# Row 65
$ws.Range("H65").Value = 8356.429
$ws.Range("I65").Value = 9339
$ws.Range("J65").Value = 5900
$ws.Range("K65").Value = 46695
$ws.Range("L65").Value = 29500
$ws.Range("M65").Value = -43575
$ws.Range("N65").Value = -35740
# Row 99
$ws.Range("H99").Value = 850
$ws.Range("I99").Value = 850
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 850
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 648
$ws.Range("N99").Value = $null
# Row 126
$ws.Range("H126").Value = 850
$ws.Range("I126").Value = 850
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 2550
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -80
$ws.Range("N126").Value = $null

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1998.5883
$ws.Range("I5").Value = 990
$ws.Range("J5").Value = 2418.8333
$ws.Range("K5").Value = 2970
$ws.Range("L5").Value = 7256.499899999999
$ws.Range("M5").Value = -2858
$ws.Range("N5").Value = -7480.499899999999
# Row 131
$ws.Range("H131").Value = 1564.1428
$ws.Range("I131").Value = 1433
$ws.Range("J131").Value = 1662.5
$ws.Range("K131").Value = 4299
$ws.Range("L131").Value = 4987.5
$ws.Range("M131").Value = 741
$ws.Range("N131").Value = -15067.5
# Row 135
$ws.Range("H135").Value = 1998.5883
$ws.Range("I135").Value = 990
$ws.Range("J135").Value = 2418.8333
$ws.Range("K135").Value = 8910
$ws.Range("L135").Value = 21769.4997
$ws.Range("M135").Value = -6375
$ws.Range("N135").Value = -26839.4997
# Row 140
$ws.Range("H140").Value = 358.33334
$ws.Range("I140").Value = 358.33334
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 1075.00002
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = 4104.999980000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 8988
$ws.Range("I80").Value = 8976
$ws.Range("J80").Value = 9000
$ws.Range("K80").Value = 8976
$ws.Range("L80").Value = 9000
$ws.Range("M80").Value = -7978
$ws.Range("N80").Value = -10996
# Row 83
$ws.Range("H83").Value = 8988
$ws.Range("I83").Value = 8976
$ws.Range("J83").Value = 9000
$ws.Range("K83").Value = 44880
$ws.Range("L83").Value = 45000
$ws.Range("M83").Value = -39888
$ws.Range("N83").Value = -54984
# Row 102
$ws.Range("H102").Value = 5174.5
$ws.Range("I102").Value = 6750
$ws.Range("J102").Value = 3599
$ws.Range("K102").Value = 6750
$ws.Range("L102").Value = 3599
$ws.Range("M102").Value = -5128
$ws.Range("N102").Value = -6843
# Row 126
$ws.Range("H126").Value = 6662.6665
$ws.Range("I126").Value = 6662.6665
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 19987.9995
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -17517.9995
# Row 132
$ws.Range("H132").Value = 6355
$ws.Range("I132").Value = 6888.5
$ws.Range("J132").Value = 5999.3335
$ws.Range("K132").Value = 20665.5
$ws.Range("L132").Value = 17998.0005
$ws.Range("M132").Value = -18135.5
$ws.Range("N132").Value = -23058.0005

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4563.2144
$ws.Range("I7").Value = 3820.5557
$ws.Range("J7").Value = 5900
$ws.Range("K7").Value = 3820.5557
$ws.Range("L7").Value = 5900
$ws.Range("M7").Value = -3708.5557
$ws.Range("N7").Value = -6124
# Row 46
$ws.Range("H46").Value = 816.6667
$ws.Range("I46").Value = 725
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 725
$ws.Range("L46").Value = 1000
$ws.Range("M46").Value = -537
$ws.Range("N46").Value = -1376
# Row 55
$ws.Range("H55").Value = 998.3077
$ws.Range("I55").Value = 463
$ws.Range("J55").Value = 1457.1428
$ws.Range("K55").Value = 463
$ws.Range("L55").Value = 1457.1428
$ws.Range("M55").Value = -290
$ws.Range("N55").Value = -1803.1428
# Row 82
$ws.Range("H82").Value = 1736.9231
$ws.Range("I82").Value = 1768.8
$ws.Range("J82").Value = 1630.6666
$ws.Range("K82").Value = 1768.8
$ws.Range("L82").Value = 1630.6666
$ws.Range("M82").Value = -1407.8
$ws.Range("N82").Value = -2352.6666
# Row 85
$ws.Range("H85").Value = 1736.9231
$ws.Range("I85").Value = 1768.8
$ws.Range("J85").Value = 1630.6666
$ws.Range("K85").Value = 1768.8
$ws.Range("L85").Value = 1630.6666
$ws.Range("M85").Value = -520.8
$ws.Range("N85").Value = -4126.6666
# Row 126
$ws.Range("H126").Value = 4563.2144
$ws.Range("I126").Value = 3820.5557
$ws.Range("J126").Value = 5900
$ws.Range("K126").Value = 11461.6671
$ws.Range("L126").Value = 17700
$ws.Range("M126").Value = -8991.667099999999
$ws.Range("N126").Value = -22640
# Row 132
$ws.Range("H132").Value = 3209.0908
$ws.Range("I132").Value = 3209.0908
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9627.2724
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7097.2724
# Row 136
$ws.Range("H136").Value = 3258.7693
$ws.Range("I136").Value = 3258.7693
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 9776.3079
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -7226.3079

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 3973.6
$ws.Range("I132").Value = 1883.9166
$ws.Range("J132").Value = 12332.333
$ws.Range("K132").Value = 5651.7498
$ws.Range("L132").Value = 36996.999
$ws.Range("M132").Value = -3121.7498
$ws.Range("N132").Value = -42056.999
